$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in attendance row for 10/13 meeting (row 21)
$ws.Range("C20").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("B21").Value = "10/13 / 4:15"
$ws.Range("C21").Value = "Google Hangout"
$ws.Range("D21").Value = "A"
$ws.Range("E21").Value = "A"
$ws.Range("F21").Value = "A"
$ws.Range("G21").Value = "A"
$ws.Range("H21").Value = "U"
$ws.Range("I21").Value = "A"

# Update selection / view
$ws.Range("D21").Select()
